$d = $word.ActiveDocument

# Find the paragraph containing the leftover author comment and delete the
# whole paragraph (including its paragraph mark).
$r = $d.Content
$found = $r.Find.Execute(
    "(when we move this over to real doc we should try to lump all equations together to save space.)",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $p = $r.Paragraphs(1)
    $p.Range.Delete()
}
